# "Add hall to pin table"
#
# 1. Fix the Wheel Motor1 "Enable" pin on the Usages sheet (row 2): it was
#    mistakenly recorded as D12, correct it to D2.
# 2. Insert a new "Hall Sensor" entry (Category=Hall Sensor, Function=<blank>,
#    Pin=A3) right above the "Color Sensor" rows on the Usages sheet, styled
#    like the other "Category header" rows (a solid accent fill, based on the
#    "60% - Accent4" cell style family).
# 3. Leave the "Used"/"Free" helper formulas on "Arduino Mega Pins" alone -
#    they recalc automatically off of the Usages!$C:$C column.

$wb = $excel.ActiveWorkbook
$usages = $wb.Worksheets.Item("Usages")
$pins = $wb.Worksheets.Item("Arduino Mega Pins")

# --- 1. Correct Wheel Motor1 / Enable pin: D12 -> D2 ---------------------
$usages.Range("C2").Value = "D2"

# --- 2. Insert the new Hall Sensor row at row 52 --------------------------
$usages.Rows.Item(52).Insert()

$usages.Range("A52:C52").Style = "60% - Accent4"
$usages.Range("A52").Value = "Hall Sensor"
$usages.Range("B52").Value = ""
$usages.Range("C52").Value = "A3"

$usages.Range("A52:C52").Interior.ThemeColor = 9

# --- 3. Restore the selections roughly where the author left them --------
$pins.Activate()
$pins.Range("B14").Select()

$usages.Activate()
$usages.Range("F37").Select()
